# Append three new incident rows (36-38) to the "Optical_Power" sheet,
# extending the used range from A1:N35 to A1:N38.
#
# Columns:
#   A Caso (case #)            - text, numeric-looking -> needs leading '
#   B F. De Reclamo (date)     - text "M/D/YYYY"        -> needs leading '
#   C Direccion                - text
#   D Comuna                   - number
#   E OT                       - text
#   F Proveedor Asignado       - text
#   G Estado                   - text
#   H Observaciones            - text
#   I Attachments              - number
#   J API_Response (JSON)      - text
#   K Coordenada_X             - number
#   L Coordenada_Y             - number
#   M Operacion                - text
#   N Zona                     - text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optical_Power")

# ---- Row 36 ----
$ws.Cells.Item(36, 1).Value  = "'5787"
$ws.Cells.Item(36, 2).Value  = "'9/22/2025"
$ws.Cells.Item(36, 3).Value  = "GAONA AV. 4102"
$ws.Cells.Item(36, 4).Value  = 10
$ws.Cells.Item(36, 5).Value  = "Pendiente ADM"
$ws.Cells.Item(36, 6).Value  = "Optical Power"
$ws.Cells.Item(36, 7).Value  = "Pendiente"
$ws.Cells.Item(36, 8).Value  = "Bajada de cliente colgando a baja altura"
$ws.Cells.Item(36, 9).Value  = 1
$ws.Cells.Item(36, 10).Value = '{"direccionesNormalizadas": [{"altura": 4102, "cod_calle": 7025, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.481890", "y": "-34.624370"}, "direccion": "GAONA AV. 4102, CABA", "nombre_calle": "GAONA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(36, 11).Value = -58.481890
$ws.Cells.Item(36, 12).Value = -34.624370
$ws.Cells.Item(36, 13).Value = "Devoto"
$ws.Cells.Item(36, 14).Value = "Capital Norte"

# ---- Row 37 ----
$ws.Cells.Item(37, 1).Value  = "'3496"
$ws.Cells.Item(37, 2).Value  = "'9/22/2025"
$ws.Cells.Item(37, 3).Value  = "CASEROS AV. 3017"
$ws.Cells.Item(37, 4).Value  = 4
$ws.Cells.Item(37, 5).Value  = "Pendiente ADM"
$ws.Cells.Item(37, 6).Value  = "Optical Power"
$ws.Cells.Item(37, 7).Value  = "Pendiente"
$ws.Cells.Item(37, 8).Value  = "Tendido a baja altura"
$ws.Cells.Item(37, 9).Value  = 1
$ws.Cells.Item(37, 10).Value = '{"direccionesNormalizadas": [{"altura": 3017, "cod_calle": 3085, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.405983", "y": "-34.637104"}, "direccion": "CASEROS AV. 3017, CABA", "nombre_calle": "CASEROS AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(37, 11).Value = -58.405983
$ws.Cells.Item(37, 12).Value = -34.637104
$ws.Cells.Item(37, 13).Value = "San Telmo"
$ws.Cells.Item(37, 14).Value = "Capital Sur"

# ---- Row 38 ----
$ws.Cells.Item(38, 1).Value  = "'3498"
$ws.Cells.Item(38, 2).Value  = "'9/22/2025"
$ws.Cells.Item(38, 3).Value  = "CASEROS AV. 3219"
$ws.Cells.Item(38, 4).Value  = 4
$ws.Cells.Item(38, 5).Value  = "Pendiente ADM"
$ws.Cells.Item(38, 6).Value  = "Optical Power"
$ws.Cells.Item(38, 7).Value  = "Pendiente"
$ws.Cells.Item(38, 8).Value  = "Tendido a baja altura"
$ws.Cells.Item(38, 9).Value  = 1
$ws.Cells.Item(38, 10).Value = '{"direccionesNormalizadas": [{"altura": 3219, "cod_calle": 3085, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.409480", "y": "-34.637709"}, "direccion": "CASEROS AV. 3219, CABA", "nombre_calle": "CASEROS AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(38, 11).Value = -58.409480
$ws.Cells.Item(38, 12).Value = -34.637709
$ws.Cells.Item(38, 13).Value = "San Telmo"
$ws.Cells.Item(38, 14).Value = "Capital Sur"
